# Add two new columns, I ("I0") and J ("IF"), to the data table on Sheet1.
# Header row (row 1) gets the new labels with the same formatting as the
# existing header cells (bold font, thin border, centered/top aligned) —
# achieved by copying H1's format onto I1:J1 rather than rebuilding the
# style by hand.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-74) ----------------------------------------------------
$values = @(
    @(7,8),
    @(6,7),
    @(8,8),
    @(7,7),
    @(7,8),
    @(9,9),
    @(8,8),
    @(5,6),
    @(7,8),
    @(8,8),
    @(6,6),
    @(8,8),
    @(8,9),
    @(8,8),
    @(8,8),
    @(7,8),
    @(7,8),
    @(9,9),
    @(7,7),
    @(9,9),
    @(8,9),
    @(6,7),
    @(8,8),
    @(10,10),
    @(6,7),
    @(8,8),
    @(9,9),
    @(6,6),
    @(9,9),
    @(7,7),
    @(6,7),
    @(9,9),
    @(7,8),
    @(9,9),
    @(6,7),
    @(6,8),
    @(5,7),
    @(8,8),
    @(6,7),
    @(8,9),
    @(7,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(6,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(6,7),
    @(8,8),
    @(9,9),
    @(7,8),
    @(6,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(7,9),
    @(9,9),
    @(6,8),
    @(10,11),
    @(8,8),
    @(8,9),
    @(6,7),
    @(5,5),
    @(8,8),
    @(3,4),
    @(7,7),
    @(8,8),
    @(7,7)
)

for ($k = 0; $k -lt $values.Count; $k++) {
    $row = 2 + $k
    $ws.Cells.Item($row, 9).Value  = $values[$k][0]   # column I
    $ws.Cells.Item($row, 10).Value = $values[$k][1]   # column J
}
